{"js": "// The page footer block (\"Ver no Jupiter Salvar em pdf Salvar em docx\" /\n// \"\u00a9 2020 . Contact: ...\") together with the blank paragraph that precedes\n// it was removed from the document (site rebuild no longer emits it).\n// Locate the footer text, then remove it plus the blank paragraph right\n// before it and the copyright paragraph right after it.\n\nconst searchResults = context.document.body.search(\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  { matchCase: true, matchWholeWord: false }\n);\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Target paragraph 'Ver no Jupiter...' not found\");\n}\n\nconst jupiterParagraph = searchResults.items[0].paragraphs.getFirst();\nconst blankParagraph = jupiterParagraph.getPrevious();\nconst copyrightParagraph = jupiterParagraph.getNext();\n\n// Confirm we are removing the expected sibling paragraphs before mutating.\nblankParagraph.load(\"text\");\ncopyrightParagraph.load(\"text\");\nawait context.sync();\n\nif (blankParagraph.text.trim() !== \"\") {\n  throw new Error(\"Paragraph before 'Ver no Jupiter...' was not blank as expected\");\n}\nif (!copyrightParagraph.text.startsWith(\"\u00a9 2020\")) {\n  throw new Error(\"Paragraph after 'Ver no Jupiter...' was not the copyright line as expected\");\n}\n\nblankParagraph.delete();\njupiterParagraph.delete();\ncopyrightParagraph.delete();\n\nawait context.sync();\n", "ps1": "# The page footer block (\"Ver no Jupiter Salvar em pdf Salvar em docx\" /\n# \"\u00a9 2020 . Contact: ...\") together with the blank paragraph that precedes\n# it was removed from the document (site rebuild no longer emits it).\n# Locate the footer text, then remove it plus the blank paragraph right\n# before it and the copyright paragraph right after it.\n\n$d = $word.ActiveDocument\n\n$findRange = $d.Content\n$findRange.Find.ClearFormatting()\n$found = $findRange.Find.Execute(\"Ver no Jupiter Salvar em pdf Salvar em docx\")\nif (-not $found) {\n    throw \"Target paragraph 'Ver no Jupiter...' not found\"\n}\n\n# $findRange now spans just the matched text; resolve which paragraph (by\n# 1-based index, Word convention) in the document it belongs to.\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $candidate = $d.Paragraphs.Item($i)\n    if ($candidate.Range.Start -le $findRange.Start -and $candidate.Range.End -ge $findRange.End) {\n        $targetIndex = $i\n        break\n    }\n}\nif ($targetIndex -eq -1) {\n    throw \"Could not resolve paragraph index for the matched text\"\n}\n\n$blankParagraph = $d.Paragraphs.Item($targetIndex - 1)\n$jupiterParagraph = $d.Paragraphs.Item($targetIndex)\n$copyrightParagraph = $d.Paragraphs.Item($targetIndex + 1)\n\nif ($blankParagraph.Range.Text.Trim() -ne \"\") {\n    throw \"Paragraph before 'Ver no Jupiter...' was not blank as expected\"\n}\nif (-not $copyrightParagraph.Range.Text.StartsWith(\"\u00a9 2020\")) {\n    throw \"Paragraph after 'Ver no Jupiter...' was not the copyright line as expected\"\n}\n\n# Delete from the bottom up so earlier indices stay valid.\n$copyrightParagraph.Range.Delete()\n$jupiterParagraph.Range.Delete()\n$blankParagraph.Range.Delete()\n"}
